$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Describe Cuarentena" text (column O, rows 10-307) - all rows
# share the same description string, which changed the quarantine hours
# from 6:00am-8:00pm to 5:00am-9:00pm and dropped the "circulando dos
# digitos por dia" clause.
$newText = "Segmentando a la poblacion para poder circular conforme a la terminacion de los digitos de su tarjeta de identidad, pasaporte o carnet de residente para extranjeros, para que puedan abastecerse de insumos básicos,  con horario de 5:00 am a 9:00 pm.  De lunes a domingo."
$ws.Range("O10:O307").Value = $newText

# Update the selection / scroll position saved in the sheet view.
$ws.Range("K12").Select()
